$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore column A amino-acid labels (row header) to guard against round-trip
# string corruption, then apply the updated PSSM values for columns B:K.
$ws.Range("A2").Value = "F"
$ws.Range("A3").Value = "W"
$ws.Range("A4").Value = "Y"
$ws.Range("A5").Value = "P"
$ws.Range("A6").Value = "M"
$ws.Range("A7").Value = "I"
$ws.Range("A8").Value = "L"
$ws.Range("A9").Value = "V"
$ws.Range("A10").Value = "A"
$ws.Range("A11").Value = "G"
$ws.Range("A12").Value = "C"
$ws.Range("A13").Value = "S"
$ws.Range("A14").Value = "T"
$ws.Range("A15").Value = "N"
$ws.Range("A16").Value = "Q"
$ws.Range("A17").Value = "D"
$ws.Range("A18").Value = "E"
$ws.Range("A19").Value = "H"
$ws.Range("A20").Value = "K"
$ws.Range("A21").Value = "R"

$ws.Range("B2").Value = -18.76145708749627
$ws.Range("C2").Value = 2.535625420575654
$ws.Range("D2").Value = -18.76145708749627
$ws.Range("E2").Value = -18.76145708749627
$ws.Range("F2").Value = -18.76145708749627
$ws.Range("G2").Value = -18.76145708749627
$ws.Range("H2").Value = -18.76145708749627
$ws.Range("I2").Value = -18.76145708749627
$ws.Range("J2").Value = -18.76145708749627
$ws.Range("K2").Value = -18.76145708749627
$ws.Range("B3").Value = -18.76145708749627
$ws.Range("C3").Value = -18.76145708749627
$ws.Range("D3").Value = -18.76145708749627
$ws.Range("E3").Value = -18.76145708749627
$ws.Range("F3").Value = -18.76145708749627
$ws.Range("G3").Value = -18.76145708749627
$ws.Range("H3").Value = -18.76145708749627
$ws.Range("I3").Value = 2.367650437506193
$ws.Range("J3").Value = -18.76145708749627
$ws.Range("K3").Value = -18.76145708749627
$ws.Range("B4").Value = -18.76145708749627
$ws.Range("C4").Value = 1.873018095764487
$ws.Range("D4").Value = -18.76145708749627
$ws.Range("E4").Value = -18.76145708749627
$ws.Range("F4").Value = 2.67527122694673
$ws.Range("G4").Value = -18.76145708749627
$ws.Range("H4").Value = 1.823410934328471
$ws.Range("I4").Value = -18.76145708749627
$ws.Range("J4").Value = 2.435596986999152
$ws.Range("K4").Value = -18.76145708749627
$ws.Range("B5").Value = -18.76145708749627
$ws.Range("C5").Value = 1.255238322033517
$ws.Range("D5").Value = -18.76145708749627
$ws.Range("E5").Value = -18.76145708749627
$ws.Range("F5").Value = -18.76145708749627
$ws.Range("G5").Value = 2.107079895578088
$ws.Range("H5").Value = -18.76145708749627
$ws.Range("I5").Value = -18.76145708749627
$ws.Range("J5").Value = -18.76145708749627
$ws.Range("K5").Value = -18.76145708749627
$ws.Range("B6").Value = -18.76145708749627
$ws.Range("C6").Value = -18.76145708749627
$ws.Range("D6").Value = -18.76145708749627
$ws.Range("E6").Value = -18.76145708749627
$ws.Range("F6").Value = -18.76145708749627
$ws.Range("G6").Value = -18.76145708749627
$ws.Range("H6").Value = -18.76145708749627
$ws.Range("I6").Value = -18.76145708749627
$ws.Range("J6").Value = -18.76145708749627
$ws.Range("K6").Value = -18.76145708749627
$ws.Range("B7").Value = 2.989215964146267
$ws.Range("C7").Value = -18.76145708749627
$ws.Range("D7").Value = -18.76145708749627
$ws.Range("E7").Value = -18.76145708749627
$ws.Range("F7").Value = -18.76145708749627
$ws.Range("G7").Value = -18.76145708749627
$ws.Range("H7").Value = -18.76145708749627
$ws.Range("I7").Value = -18.76145708749627
$ws.Range("J7").Value = -18.76145708749627
$ws.Range("K7").Value = -18.76145708749627
$ws.Range("B8").Value = -18.76145708749627
$ws.Range("C8").Value = -18.76145708749627
$ws.Range("D8").Value = -18.76145708749627
$ws.Range("E8").Value = 2.816440523291857
$ws.Range("F8").Value = -18.76145708749627
$ws.Range("G8").Value = -18.76145708749627
$ws.Range("H8").Value = -18.76145708749627
$ws.Range("I8").Value = -18.76145708749627
$ws.Range("J8").Value = -18.76145708749627
$ws.Range("K8").Value = -18.76145708749627
$ws.Range("B9").Value = 3.592102487805461
$ws.Range("C9").Value = -18.76145708749627
$ws.Range("D9").Value = -18.76145708749627
$ws.Range("E9").Value = -18.76145708749627
$ws.Range("F9").Value = -18.76145708749627
$ws.Range("G9").Value = -18.76145708749627
$ws.Range("H9").Value = -18.76145708749627
$ws.Range("I9").Value = -18.76145708749627
$ws.Range("J9").Value = -18.76145708749627
$ws.Range("K9").Value = -18.76145708749627
$ws.Range("B10").Value = -18.76145708749627
$ws.Range("C10").Value = -18.76145708749627
$ws.Range("D10").Value = -18.76145708749627
$ws.Range("E10").Value = -18.76145708749627
$ws.Range("F10").Value = -18.76145708749627
$ws.Range("G10").Value = -18.76145708749627
$ws.Range("H10").Value = -18.76145708749627
$ws.Range("I10").Value = 1.566732427953247
$ws.Range("J10").Value = -18.76145708749627
$ws.Range("K10").Value = 2.25516963883328
$ws.Range("B11").Value = -18.76145708749627
$ws.Range("C11").Value = -18.76145708749627
$ws.Range("D11").Value = -18.76145708749627
$ws.Range("E11").Value = 1.953946100014616
$ws.Range("F11").Value = -18.76145708749627
$ws.Range("G11").Value = 2.530824660139925
$ws.Range("H11").Value = -18.76145708749627
$ws.Range("I11").Value = -18.76145708749627
$ws.Range("J11").Value = -18.76145708749627
$ws.Range("K11").Value = 1.356587952095764
$ws.Range("B12").Value = -18.76145708749627
$ws.Range("C12").Value = -18.76145708749627
$ws.Range("D12").Value = -18.76145708749627
$ws.Range("E12").Value = -18.76145708749627
$ws.Range("F12").Value = -18.76145708749627
$ws.Range("G12").Value = -18.76145708749627
$ws.Range("H12").Value = -18.76145708749627
$ws.Range("I12").Value = -18.76145708749627
$ws.Range("J12").Value = -18.76145708749627
$ws.Range("K12").Value = -18.76145708749627
$ws.Range("B13").Value = -18.76145708749627
$ws.Range("C13").Value = -18.76145708749627
$ws.Range("D13").Value = -18.76145708749627
$ws.Range("E13").Value = 1.729790258153019
$ws.Range("F13").Value = -18.76145708749627
$ws.Range("G13").Value = -18.76145708749627
$ws.Range("H13").Value = -18.76145708749627
$ws.Range("I13").Value = -18.76145708749627
$ws.Range("J13").Value = 2.235505254894622
$ws.Range("K13").Value = 1.576671008582199
$ws.Range("B14").Value = -18.76145708749627
$ws.Range("C14").Value = -18.76145708749627
$ws.Range("D14").Value = 4.321925010724788
$ws.Range("E14").Value = -18.76145708749627
$ws.Range("F14").Value = -18.76145708749627
$ws.Range("G14").Value = -18.76145708749627
$ws.Range("H14").Value = -18.76145708749627
$ws.Range("I14").Value = -18.76145708749627
$ws.Range("J14").Value = -18.76145708749627
$ws.Range("K14").Value = 2.073736739958098
$ws.Range("B15").Value = -18.76145708749627
$ws.Range("C15").Value = -18.76145708749627
$ws.Range("D15").Value = -18.76145708749627
$ws.Range("E15").Value = -18.76145708749627
$ws.Range("F15").Value = -18.76145708749627
$ws.Range("G15").Value = -18.76145708749627
$ws.Range("H15").Value = -18.76145708749627
$ws.Range("I15").Value = -18.76145708749627
$ws.Range("J15").Value = -18.76145708749627
$ws.Range("K15").Value = -18.76145708749627
$ws.Range("B16").Value = -18.76145708749627
$ws.Range("C16").Value = -18.76145708749627
$ws.Range("D16").Value = -18.76145708749627
$ws.Range("E16").Value = -18.76145708749627
$ws.Range("F16").Value = -18.76145708749627
$ws.Range("G16").Value = -18.76145708749627
$ws.Range("H16").Value = -18.76145708749627
$ws.Range("I16").Value = -18.76145708749627
$ws.Range("J16").Value = 2.318070014250524
$ws.Range("K16").Value = -18.76145708749627
$ws.Range("B17").Value = -18.76145708749627
$ws.Range("C17").Value = 0.068516437728878
$ws.Range("D17").Value = -18.76145708749627
$ws.Range("E17").Value = -18.76145708749627
$ws.Range("F17").Value = -18.76145708749627
$ws.Range("G17").Value = -18.76145708749627
$ws.Range("H17").Value = 0.4675517895164172
$ws.Range("I17").Value = 0.8916346375015891
$ws.Range("J17").Value = 1.251725168033324
$ws.Range("K17").Value = -18.76145708749627
$ws.Range("B18").Value = -18.76145708749627
$ws.Range("C18").Value = -18.76145708749627
$ws.Range("D18").Value = -18.76145708749627
$ws.Range("E18").Value = -18.76145708749627
$ws.Range("F18").Value = -18.76145708749627
$ws.Range("G18").Value = -18.76145708749627
$ws.Range("H18").Value = 0.3308277684660155
$ws.Range("I18").Value = 0.8738518588970628
$ws.Range("J18").Value = 1.329354433484573
$ws.Range("K18").Value = -18.76145708749627
$ws.Range("B19").Value = -18.76145708749627
$ws.Range("C19").Value = -18.76145708749627
$ws.Range("D19").Value = -18.76145708749627
$ws.Range("E19").Value = -18.76145708749627
$ws.Range("F19").Value = -18.76145708749627
$ws.Range("G19").Value = -18.76145708749627
$ws.Range("H19").Value = 1.858194239745315
$ws.Range("I19").Value = 2.096126845872773
$ws.Range("J19").Value = -18.76145708749627
$ws.Range("K19").Value = -18.76145708749627
$ws.Range("B20").Value = -18.76145708749627
$ws.Range("C20").Value = 1.584053198857955
$ws.Range("D20").Value = -18.76145708749627
$ws.Range("E20").Value = -18.76145708749627
$ws.Range("F20").Value = 3.766846540731454
$ws.Range("G20").Value = -18.76145708749627
$ws.Range("H20").Value = 2.237709587918229
$ws.Range("I20").Value = 1.968353651760773
$ws.Range("J20").Value = -18.76145708749627
$ws.Range("K20").Value = 2.4522553038227
$ws.Range("B21").Value = -18.76145708749627
$ws.Range("C21").Value = 2.037344548858551
$ws.Range("D21").Value = -18.76145708749627
$ws.Range("E21").Value = 2.527224840475989
$ws.Range("F21").Value = -18.76145708749627
$ws.Range("G21").Value = 3.309286579413698
$ws.Range("H21").Value = 2.4537257229103
$ws.Range("I21").Value = -18.76145708749627
$ws.Range("J21").Value = -18.76145708749627
$ws.Range("K21").Value = -18.76145708749627
